# Apply updated crypto price/volume figures from the latest GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.946.65'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '''1.857.72'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''311.96'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D7').Value = '''0.5138'
$ws.Range('E7').Value = '  +1.70%  '
$ws.Range('D8').Value = '''0.3825'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = '''0.08233'
$ws.Range('E9').Value = '  -4.59%  '
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').Value = '''41.43'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('D14').Value = '''1.861.33'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '''7.245'
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').Value = '''1.003'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '''90.32'
$ws.Range('D19').Value = '''0.06640'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '''17.65'
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('D21').Value = '''1.002'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '''6.002'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').Value = '''27.983.12'
$ws.Range('E24').Value = '  -3.53%  '
$ws.Range('D25').Value = '''2.238'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').Value = '''2.072.17'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = '''2.501'
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').Value = '''157.26'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').Value = '''124.48'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('D31').Value = '''0.1063'
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').Value = '''1.025'
$ws.Range('E32').Value = '  -3.27%  '
$ws.Range('D33').Value = '''5.862'
$ws.Range('E33').Value = '  +4.86%  '
$ws.Range('D34').Value = '''3.590'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '''9.395'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('D36').Value = '''0.02412'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('D37').Value = '''0.06497'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('D39').Value = '''0.6534'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D41').Value = '''4.972'
$ws.Range('E41').Value = '  +1.58%  '
$ws.Range('D42').Value = '''1.205'
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('D44').Value = '''0.6109'
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('D45').Value = '''13.01'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '''3.674'
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '''1.273'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').Value = '''2.006'
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('D49').Value = '''1.211'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('D50').Value = '''120.66'
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').Value = '''78.13'
$ws.Range('E51').Value = '  -2.62%  '
